$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 46
$ws.Range("I2").Value = 128
$ws.Range("J2").Value = 580
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 165
$ws.Range("N2").Value = 89
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("R2").Value = 4
$ws.Range("T2").Value = 95
$ws.Range("V2").Value = 894
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 875
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 11
$ws.Range("AA2").Value = 6
